$wb = $excel.ActiveWorkbook
$profile = $wb.Worksheets.Item("Profile")

# Add the new "Product" worksheet right after "Profile"
$newSheet = $wb.Worksheets.Add($null, $profile)
$newSheet.Name = "Product"

# Header row (row 3)
$newSheet.Range("A3").Value = "ID"
$newSheet.Range("B3").Value = "NAME"
$newSheet.Range("C3").Value = "QUANTITY"
$newSheet.Range("A3:C3").Font.Bold = $true

# Data rows
$newSheet.Range("A4").Value = "A001"
$newSheet.Range("B4").Value = "Cheese"
$newSheet.Range("C4").Value = 24

$newSheet.Range("A5").Value = "A002"
$newSheet.Range("B5").Value = "Butter"
$newSheet.Range("C5").Value = 17

$newSheet.Range("A6").Value = "A003"
$newSheet.Range("B6").Value = "Milk"
$newSheet.Range("C6").Value = 37

# Reserved rows at the top (added last so the shared strings are appended afterwards)
$newSheet.Range("A1").Value = "Reserved row"
$newSheet.Range("A2").Value = "Reserved row"
$newSheet.Range("A1:A2").Font.Italic = $true

# Selection / active cell on the new sheet
$null = $newSheet.Range("B9").Select()

# Window / view state
$wb.Windows.Item(1).Left = 54400
$wb.Windows.Item(1).Top = 3840
$wb.Windows.Item(1).Width = 25600
$wb.Windows.Item(1).Height = 14160
